$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact text value, even when the text looks like a number
# (mirrors how Excel stores a value typed with a leading apostrophe as text).
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '67.158.56'
$ws.Range("E2").Value = '  -0.65%  '
$ws.Range("D3").Value = '2.473.40'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("E4").Value = '  -0.03%  '
Set-TextValue $ws.Range("D5") '582.46'
$ws.Range("E5").Value = '  -1.33%  '
Set-TextValue $ws.Range("D6") '168.20'
$ws.Range("E6").Value = '  -3.10%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  -1.73%  '
$ws.Range("D9").Value = '2.473.35'
$ws.Range("E9").Value = '  -0.69%  '
$ws.Range("E10").Value = '  -3.35%  '
$ws.Range("E11").Value = '  -0.98%  '
Set-TextValue $ws.Range("D12") '4.97'
$ws.Range("E12").Value = '  -2.37%  '
$ws.Range("E13").Value = '  -2.34%  '
$ws.Range("D14").Value = '2.923.76'
$ws.Range("E14").Value = '  -1.00%  '
Set-TextValue $ws.Range("D15") '25.50'
$ws.Range("E15").Value = '  -2.97%  '
$ws.Range("D16").Value = '67.115.92'
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("E17").Value = '  -3.83%  '
$ws.Range("D18").Value = '2.493.77'
$ws.Range("E18").Value = '  +1.76%  '
$ws.Range("E19").Value = '  -4.40%  '
$ws.Range("E20").Value = '  -4.68%  '
Set-TextValue $ws.Range("D21") '357.04'
$ws.Range("E21").Value = '  -2.69%  '
$ws.Range("E22").Value = '  -1.71%  '
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("E24").Value = '  -2.66%  '
$ws.Range("E25").Value = '  -6.72%  '
Set-TextValue $ws.Range("D26") '1.79'
$ws.Range("E26").Value = '  -6.98%  '
Set-TextValue $ws.Range("D27") '9.13'
$ws.Range("E27").Value = '  -8.63%  '
Set-TextValue $ws.Range("D28") '0.999'
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("D29").Value = '2.598.65'
$ws.Range("E29").Value = '  -0.92%  '
$ws.Range("D30").Value = '0.0₃0906'
$ws.Range("E30").Value = '  -5.33%  '
Set-TextValue $ws.Range("D31") '510.10'
$ws.Range("E31").Value = '  -4.08%  '
Set-TextValue $ws.Range("D32") '7.78'
$ws.Range("E32").Value = '  -6.46%  '
$ws.Range("E33").Value = '  -4.21%  '
$ws.Range("E34").Value = '  -5.57%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  -6.57%  '
Set-TextValue $ws.Range("D37") '158.25'
$ws.Range("E37").Value = '  +0.35%  '
Set-TextValue $ws.Range("D39") '18.46'
$ws.Range("E39").Value = '  -1.26%  '
$ws.Range("E40").Value = '  -5.29%  '
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("E42").Value = '  -6.05%  '
$ws.Range("E43").Value = '  -5.97%  '
$ws.Range("E44").Value = '  -6.38%  '
$ws.Range("E45").Value = '  -5.99%  '
Set-TextValue $ws.Range("D47") '141.34'
$ws.Range("E47").Value = '  -2.30%  '
$ws.Range("E48").Value = '  -5.48%  '
$ws.Range("E49").Value = '  -5.65%  '
Set-TextValue $ws.Range("D50") '1.59'
$ws.Range("E50").Value = '  -5.42%  '
$ws.Range("E51").Value = '  -8.67%  '
